# Swap data between row 59 and row 60, and between row 66 and row 67,
# as described by the commit diff. Columns that are identical between
# the paired rows are left untouched; only the columns that actually
# differ are swapped.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellValue($ws, $addr, $val, $forceText) {
    if ($forceText -and $val -ne "") {
        # Numeric-looking text (e.g. "50") needs to be forced to Text so
        # Excel doesn't silently convert it to a number.
        $ws.Range($addr).NumberFormat = "@"
        $ws.Range($addr).Value = $val
        $ws.Range($addr).Style = "Normal"
    } else {
        $ws.Range($addr).Value = $val
    }
}

function Swap-Cell($ws, $col, $row1, $row2, $forceText) {
    $addr1 = "$col$row1"
    $addr2 = "$col$row2"

    $v1 = $ws.Range($addr1).Value2
    $v2 = $ws.Range($addr2).Value2

    Set-CellValue $ws $addr1 $v2 $forceText
    Set-CellValue $ws $addr2 $v1 $forceText
}

# --- Swap rows 59 and 60 ---
Swap-Cell $ws "A"  59 60 $false
Swap-Cell $ws "I"  59 60 $true
Swap-Cell $ws "J"  59 60 $false
Swap-Cell $ws "Q"  59 60 $false
Swap-Cell $ws "R"  59 60 $false
Swap-Cell $ws "Z"  59 60 $false
Swap-Cell $ws "AB" 59 60 $false
Swap-Cell $ws "AC" 59 60 $false

# --- Swap rows 66 and 67 ---
Swap-Cell $ws "A"  66 67 $false
Swap-Cell $ws "B"  66 67 $false
Swap-Cell $ws "E"  66 67 $false
Swap-Cell $ws "F"  66 67 $false
Swap-Cell $ws "G"  66 67 $false
Swap-Cell $ws "H"  66 67 $false
Swap-Cell $ws "I"  66 67 $true
Swap-Cell $ws "J"  66 67 $false
Swap-Cell $ws "Q"  66 67 $false
Swap-Cell $ws "R"  66 67 $false
Swap-Cell $ws "Z"  66 67 $false
Swap-Cell $ws "AB" 66 67 $false
Swap-Cell $ws "AC" 66 67 $false
